# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds the "K" (strikeouts) values for rows 2-16.
# New strikeout counts replacing the old ones.
$kValues = @{
    2  = 4
    3  = 0
    4  = 0
    5  = 2
    6  = 3
    7  = 5
    8  = 3
    9  = 5
    10 = 4
    11 = 5
    12 = 8
    13 = 1
    14 = 0
    15 = 1
    16 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
